$d = $word.ActiveDocument

$d.Content.Find.Execute("If you make an new agreement", $true, $false, $false, $false, $false,
                         $true, 1, $false, "If you make a new agreement", 2)

$d.Content.Find.Execute("Until March 31, 2021", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Until June 30, 2021", 2)
